$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The break log is rebuilt for a new day (2025-01-21) with fewer, non-
# overlapping entries: breaks now fit in a single line, and any breaks
# that overlap move into one extra column/layer (D) instead of spilling
# into extra rows. Only 3 rows of data remain (was 7).

# Row 1: AS break, with an overlap pushed into column D
$ws.Range("A1").Value = "AS"
$ws.Range("B1").Value = "2025-01-21T11:50"
$ws.Range("C1").Value = "2025-01-21T12:30"
$ws.Range("D1").Value = "2025-01-21T14:50"

# Row 2: EM break, no overlap layer needed
$ws.Range("A2").Value = "EM"
$ws.Range("B2").Value = "2025-01-21T12:10"
$ws.Range("C2").Value = "2025-01-21T14:30"
$ws.Range("D2").Value = ""

# Row 3: ABC break, single timestamp only
$ws.Range("A3").Value = "ABC"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "2025-01-21T12:20"
$ws.Range("D3").Value = ""

# Drop the old rows 4-7 (KG/CS/DC/AL/MD entries no longer present) so the
# sheet's used range shrinks back down to A1:D3.
$ws.Range("A4:D7").Clear()
